$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.02"
$ws.Range("D3").Value = "'21.07"
$ws.Range("D4").Value = "'6.225"
$ws.Range("D5").Value = "'0.06196"
$ws.Range("D6").Value = "'3.575"
$ws.Range("D7").Value = "'6.562"
$ws.Range("D9").Value = "'0.8246"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1620"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08208"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03507"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03096"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09129"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.774"
$ws.Range("E15").Value = "14MCDexMCBBestin24h"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001611"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04698"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006456"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.006172"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001067"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.0001502"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.808"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "'2.283"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "'0.01384"
$ws.Range("E24").Value = "23OneONE"
$ws.Range("D25").Value = "'0.3387"
$ws.Range("D26").Value = "'0.1223"
$ws.Range("D40").Value = "'0.04672"
$ws.Range("D41").Value = "'0.007050"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1106"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003525"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.01122"
$ws.Range("D45").Value = "'0.00006376"
$ws.Range("D47").Value = "'0.8453"
$ws.Range("D48").Value = "'0.002322"
